# Rename "Government" to "Public Administration" (official term used).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Cells C15 and E15 hold the English label "Government" -> rename to "Public Administration"
$ws.Range("C15").Value = "Public Administration"
$ws.Range("E15").Value = "Public Administration"

# Reflect the active cell/selection left by the edit (E15), matching the saved view state.
$ws.Range("E15").Select()
